# Applies the "agregar agrupación por user_id y session_id con conteo de
# usuarios únicos por ruta" update to the KPI results sheet:
#  - refreshes a batch of existing KPI values/descriptions (now computed
#    with the user_id/session_id grouping in place)
#  - appends four new KPI rows (29-32) describing the per-route unique
#    user grouping, distribution, concentration and diversity metrics

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Updated numeric KPI results -----------------------------------
$ws.Range("B3").Value  = 60
$ws.Range("B4").Value  = 3.6
$ws.Range("B6").Value  = 60
$ws.Range("B8").Value  = 1.054920167986144
$ws.Range("B11").Value = 16.66666666666666
$ws.Range("B12").Value = 60
$ws.Range("B13").Value = 0.4898979485566356
$ws.Range("B19").Value = 33.33333333333333
$ws.Range("B20").Value = 33.33333333333333
$ws.Range("B21").Value = 20
$ws.Range("B23").Value = 60
$ws.Range("B25").Value = 2
$ws.Range("B26").Value = 2
$ws.Range("B27").Value = 1.666666666666667

# ---- Updated text KPI results ---------------------------------------
$ws.Range("B9").Value = "1. session_start: 2 veces (40.0%)`n2. start_checkout: 2 veces (40.0%)`n3. first_visit: 1 veces (20.0%)"

$ws.Range("B10").Value = "1. begin_checkout: 2 veces (40.0%)`n2. purchase_GA4: 2 veces (40.0%)`n3. first_visit: 1 veces (20.0%)"

$ws.Range("B14").Value = "1. begin_checkout: 5 veces (100.0%)"

$ws.Range("B16").Value = "1. begin_checkout: 2 veces (66.7%)`n2. first_visit: 1 veces (33.3%)"

$ws.Range("B17").Value = "1. begin_checkout: 5 veces (27.8%)`n2. purchase_GA4: 5 veces (27.8%)`n3. first_visit: 4 veces (22.2%)`n4. session_start: 2 veces (11.1%)`n5. start_checkout: 2 veces (11.1%)"

$ws.Range("B18").Value = "1. session_start → first_visit: 2 veces (15.4%)`n2. first_visit → begin_checkout: 2 veces (15.4%)`n3. start_checkout → begin_checkout: 2 veces (15.4%)`n4. begin_checkout → purchase_GA4: 2 veces (15.4%)`n5. purchase_GA4 → purchase_GA4: 2 veces (15.4%)`n6. first_visit → purchase_GA4: 1 veces (7.7%)`n7. purchase_GA4 → begin_checkout: 1 veces (7.7%)`n8. begin_checkout → first_visit: 1 veces (7.7%)"

# ---- New KPI rows (29-32) for the user/session grouping feature -----
$ws.Range("A29").Value = "Rutas más populares por usuarios únicos"
$ws.Range("B29").Value = "1.  1 session_start - 1 first_visit - 1 begin_checkou... (1 usuarios únicos, 50.0%)`n2. first_visit - 3 purchase_GA4 - 1 begin_checkout@ -... (1 usuarios únicos, 50.0%)`n3. start_checkout - 1 begin_checkout@ - 1 purchase_GA... (1 usuarios únicos, 50.0%)"

$ws.Range("A30").Value = "Distribución de usuarios por ruta"
$ws.Range("B30").Value = "Distribución de usuarios por ruta:`n    - Media: 1.00 usuarios`n    - Mediana: 1.00 usuarios`n    - Máximo: 1 usuarios`n    - Mínimo: 1 usuarios`n    - Desviación estándar: 0.00 usuarios"

$ws.Range("A31").Value = "Rutas con alta concentración de usuarios"
$ws.Range("B31").Value = "Rutas que concentran más del 10% de usuarios (0 usuarios):`n-  1 session_start - 1 first_visit - 1 begin_checkou... (1 usuarios, 50.0%)`n- first_visit - 3 purchase_GA4 - 1 begin_checkout@ -... (1 usuarios, 50.0%)`n- start_checkout - 1 begin_checkout@ - 1 purchase_GA... (1 usuarios, 50.0%)"

$ws.Range("A32").Value = "Diversidad de rutas"
$ws.Range("B32").Value = "Diversidad de rutas:`n    - Total de rutas únicas: 3`n    - Total de usuarios únicos: 2`n    - Promedio de usuarios por ruta: 1.00`n    - Coeficiente de Gini (desigualdad): 0.000`n    - Ruta con más usuarios: 1 usuarios`n    - Ruta con menos usuarios: 1 usuarios"
